$d = $word.ActiveDocument

# Update the date line (unique text in the document)
$d.Content.Find.Execute("2023-10-12 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-13 Friday", 2) | Out-Null

# Update the division problems in the table, cell by cell (values repeat,
# so addressing by row/column avoids ambiguity)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "28÷4=7, 0") { throw "Unexpected content in cell (1,1): $($cell.Range.Text)" }
$cell.Range.Text = "99÷5=19, 4"

$cell = $t.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "84÷8=10, 4") { throw "Unexpected content in cell (1,2): $($cell.Range.Text)" }
$cell.Range.Text = "79÷7=11, 2"

$cell = $t.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "43÷9=4, 7") { throw "Unexpected content in cell (1,3): $($cell.Range.Text)" }
$cell.Range.Text = "85÷2=42, 1"

$cell = $t.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "92÷3=30, 2") { throw "Unexpected content in cell (1,4): $($cell.Range.Text)" }
$cell.Range.Text = "97÷8=12, 1"

$cell = $t.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "84÷5=16, 4") { throw "Unexpected content in cell (1,5): $($cell.Range.Text)" }
$cell.Range.Text = "58÷5=11, 3"

$cell = $t.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "87÷8=10, 7") { throw "Unexpected content in cell (5,1): $($cell.Range.Text)" }
$cell.Range.Text = "91÷6=15, 1"

$cell = $t.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "28÷6=4, 4") { throw "Unexpected content in cell (5,2): $($cell.Range.Text)" }
$cell.Range.Text = "79÷3=26, 1"

$cell = $t.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "75÷3=25, 0") { throw "Unexpected content in cell (5,3): $($cell.Range.Text)" }
$cell.Range.Text = "70÷5=14, 0"

$cell = $t.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "89÷3=29, 2") { throw "Unexpected content in cell (5,4): $($cell.Range.Text)" }
$cell.Range.Text = "39÷9=4, 3"

$cell = $t.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "54÷9=6, 0") { throw "Unexpected content in cell (5,5): $($cell.Range.Text)" }
$cell.Range.Text = "85÷8=10, 5"

$cell = $t.Cell(9, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "60÷2=30, 0") { throw "Unexpected content in cell (9,1): $($cell.Range.Text)" }
$cell.Range.Text = "83÷7=11, 6"

$cell = $t.Cell(9, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "99÷8=12, 3") { throw "Unexpected content in cell (9,2): $($cell.Range.Text)" }
$cell.Range.Text = "42÷2=21, 0"

$cell = $t.Cell(9, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "49÷2=24, 1") { throw "Unexpected content in cell (9,3): $($cell.Range.Text)" }
$cell.Range.Text = "69÷3=23, 0"

$cell = $t.Cell(9, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "35÷4=8, 3") { throw "Unexpected content in cell (9,4): $($cell.Range.Text)" }
$cell.Range.Text = "82÷7=11, 5"

$cell = $t.Cell(9, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "61÷9=6, 7") { throw "Unexpected content in cell (9,5): $($cell.Range.Text)" }
$cell.Range.Text = "13÷6=2, 1"

$cell = $t.Cell(13, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "37÷2=18, 1") { throw "Unexpected content in cell (13,1): $($cell.Range.Text)" }
$cell.Range.Text = "78÷5=15, 3"

$cell = $t.Cell(13, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "17÷2=8, 1") { throw "Unexpected content in cell (13,2): $($cell.Range.Text)" }
$cell.Range.Text = "85÷2=42, 1"

$cell = $t.Cell(13, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "63÷4=15, 3") { throw "Unexpected content in cell (13,3): $($cell.Range.Text)" }
$cell.Range.Text = "63÷6=10, 3"

$cell = $t.Cell(13, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "30÷9=3, 3") { throw "Unexpected content in cell (13,4): $($cell.Range.Text)" }
$cell.Range.Text = "39÷6=6, 3"

$cell = $t.Cell(13, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "81÷5=16, 1") { throw "Unexpected content in cell (13,5): $($cell.Range.Text)" }
$cell.Range.Text = "28÷5=5, 3"

$cell = $t.Cell(17, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "57÷6=9, 3") { throw "Unexpected content in cell (17,1): $($cell.Range.Text)" }
$cell.Range.Text = "67÷4=16, 3"

$cell = $t.Cell(17, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "92÷3=30, 2") { throw "Unexpected content in cell (17,2): $($cell.Range.Text)" }
$cell.Range.Text = "34÷2=17, 0"

$cell = $t.Cell(17, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "75÷8=9, 3") { throw "Unexpected content in cell (17,3): $($cell.Range.Text)" }
$cell.Range.Text = "33÷2=16, 1"

$cell = $t.Cell(17, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "92÷3=30, 2") { throw "Unexpected content in cell (17,4): $($cell.Range.Text)" }
$cell.Range.Text = "60÷4=15, 0"

$cell = $t.Cell(17, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "73÷8=9, 1") { throw "Unexpected content in cell (17,5): $($cell.Range.Text)" }
$cell.Range.Text = "70÷3=23, 1"

